$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing expense entry (Food -> Spotify subscription)
$ws.Range("A2").Value = "Spotify subscription"
$ws.Range("B2").Value = 120
$ws.Range("C2").Value = 45972.22928240741

# New row 3: Travel
$ws.Range("A3").Value = "Travel"
$ws.Range("B3").Value = 450
$ws.Range("C3").Value = 45972.22928240741

# New row 4: Groceries
$ws.Range("A4").Value = "Groceries"
$ws.Range("B4").Value = 430
$ws.Range("C4").Value = 45971.22928240741

# New row 5: Movie night
$ws.Range("A5").Value = "Movie night"
$ws.Range("B5").Value = 550
$ws.Range("C5").Value = 45970.22928240741

# Copy the date-formatted style from C2 onto the new date cells so they
# share the same number format (style index) rather than allocating a
# brand new custom number format.
$ws.Range("C2").Copy()
$ws.Range("C3:C5").PasteSpecial(-4122)
